$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper cell (outside the used range) used to force text-typed values
# into D-column cells without altering cell styles: a formula producing a
# string is computed in $helper, copied, then pasted as values-only into
# the destination cell (PasteSpecial xlPasteValues = -4163), which carries
# over the Text type without touching NumberFormat/style.
$helper = "Z1"

$ws.Range($helper).Formula = "=""70.323.51"""
$ws.Range($helper).Copy()
$ws.Range("D2").PasteSpecial(-4163)
$ws.Range($helper).Clear()
$ws.Range("E2").Value = "  +1.59%  "

$ws.Range($helper).Formula = "=""3.966.03"""
$ws.Range($helper).Copy()
$ws.Range("D3").PasteSpecial(-4163)
$ws.Range($helper).Clear()
$ws.Range("E3").Value = "  +2.49%  "

$ws.Range($helper).Formula = "=""0.997"""
$ws.Range($helper).Copy()
$ws.Range("D4").PasteSpecial(-4163)
$ws.Range($helper).Clear()
$ws.Range("E4").Value = "  -0.30%  "

$ws.Range($helper).Formula = "=""612.90"""
$ws.Range($helper).Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range($helper).Clear()
$ws.Range("E5").Value = "  +1.45%  "

$ws.Range($helper).Formula = "=""170.34"""
$ws.Range($helper).Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range($helper).Clear()
$ws.Range("E6").Value = "  +3.25%  "

$ws.Range($helper).Formula = "=""3.962.39"""
$ws.Range($helper).Copy()
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range($helper).Clear()
$ws.Range("E7").Value = "  +2.45%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range($helper).Formula = "=""0.538"""
$ws.Range($helper).Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range($helper).Clear()
$ws.Range("E9").Value = "  +0.63%  "

$ws.Range($helper).Formula = "=""0.173"""
$ws.Range($helper).Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range($helper).Clear()
$ws.Range("E10").Value = "  +1.97%  "

$ws.Range($helper).Formula = "=""6.49"""
$ws.Range($helper).Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range($helper).Clear()
$ws.Range("E11").Value = "  +2.20%  "

$ws.Range($helper).Formula = "=""0.470"""
$ws.Range($helper).Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range($helper).Clear()
$ws.Range("E12").Value = "  +1.76%  "

$ws.Range($helper).Formula = "=""0.0000258"""
$ws.Range($helper).Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range($helper).Clear()
$ws.Range("E13").Value = "  +4.80%  "

$ws.Range($helper).Formula = "=""38.24"""
$ws.Range($helper).Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range($helper).Clear()
$ws.Range("E14").Value = "  +2.72%  "

$ws.Range($helper).Formula = "=""4.598.90"""
$ws.Range($helper).Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range($helper).Clear()
$ws.Range("E15").Value = "  +1.80%  "

$ws.Range($helper).Formula = "=""3.973.24"""
$ws.Range($helper).Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range($helper).Clear()
$ws.Range("E16").Value = "  +2.61%  "

$ws.Range($helper).Formula = "=""70.008.67"""
$ws.Range($helper).Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range($helper).Clear()
$ws.Range("E17").Value = "  +0.88%  "

$ws.Range($helper).Formula = "=""7.64"""
$ws.Range($helper).Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range($helper).Clear()
$ws.Range("E18").Value = "  +0.11%  "

$ws.Range($helper).Formula = "=""17.79"""
$ws.Range($helper).Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range($helper).Clear()
$ws.Range("E19").Value = "  +3.14%  "

$ws.Range("E20").Value = "  -1.93%  "

$ws.Range($helper).Formula = "=""11.13"""
$ws.Range($helper).Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range($helper).Clear()
$ws.Range("E21").Value = "  -4.55%  "

$ws.Range($helper).Formula = "=""501.49"""
$ws.Range($helper).Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range($helper).Clear()
$ws.Range("E22").Value = "  +2.19%  "

$ws.Range($helper).Formula = "=""0.742"""
$ws.Range($helper).Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range($helper).Clear()
$ws.Range("E23").Value = "  +2.36%  "

$ws.Range($helper).Formula = "=""0.0000170"""
$ws.Range($helper).Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range($helper).Clear()
$ws.Range("E24").Value = "  +6.11%  "

$ws.Range($helper).Formula = "=""85.69"""
$ws.Range($helper).Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range($helper).Clear()
$ws.Range("E25").Value = "  +1.15%  "

$ws.Range($helper).Formula = "=""2.31"""
$ws.Range($helper).Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range($helper).Clear()
$ws.Range("E26").Value = "  +1.19%  "

$ws.Range($helper).Formula = "=""12.42"""
$ws.Range($helper).Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range($helper).Clear()
$ws.Range("E27").Value = "  +1.27%  "

$ws.Range($helper).Formula = "=""10.29"""
$ws.Range($helper).Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range($helper).Clear()
$ws.Range("E28").Value = "  +2.82%  "

$ws.Range("E29").Value = "  -0.07%  "

$ws.Range($helper).Formula = "=""3.01"""
$ws.Range($helper).Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range($helper).Clear()
$ws.Range("E30").Value = "  +0.47%  "

$ws.Range($helper).Formula = "=""4.090.98"""
$ws.Range($helper).Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range($helper).Clear()
$ws.Range("E31").Value = "  +1.71%  "

$ws.Range($helper).Formula = "=""2.43"""
$ws.Range($helper).Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range($helper).Clear()
$ws.Range("E32").Value = "  +0.90%  "

$ws.Range($helper).Formula = "=""7.92"""
$ws.Range($helper).Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range($helper).Clear()
$ws.Range("E33").Value = "  -0.88%  "

$ws.Range($helper).Formula = "=""32.51"""
$ws.Range($helper).Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range($helper).Clear()
$ws.Range("E34").Value = "  -0.15%  "

$ws.Range($helper).Formula = "=""3.907.35"""
$ws.Range($helper).Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range($helper).Clear()
$ws.Range("E35").Value = "  +2.42%  "

$ws.Range("E36").Value = "  +0.56%  "

$ws.Range($helper).Formula = "=""6.18"""
$ws.Range($helper).Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range($helper).Clear()
$ws.Range("E37").Value = "  +4.09%  "

$ws.Range($helper).Formula = "=""1.05"""
$ws.Range($helper).Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range($helper).Clear()
$ws.Range("E38").Value = "  +0.40%  "

$ws.Range($helper).Formula = "=""0.141"""
$ws.Range($helper).Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range($helper).Clear()
$ws.Range("E39").Value = "  +0.87%  "

$ws.Range("E40").Value = "  +9.00%  "

$ws.Range($helper).Formula = "=""0.996"""
$ws.Range($helper).Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range($helper).Clear()
$ws.Range("E41").Value = "  -0.44%  "

$ws.Range($helper).Formula = "=""0.327"""
$ws.Range($helper).Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range($helper).Clear()
$ws.Range("E42").Value = "  +1.99%  "

$ws.Range($helper).Formula = "=""2.07"""
$ws.Range($helper).Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range($helper).Clear()
$ws.Range("E43").Value = "  +3.71%  "

$ws.Range($helper).Formula = "=""439.11"""
$ws.Range($helper).Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range($helper).Clear()
$ws.Range("E44").Value = "  -0.52%  "

$ws.Range("E45").Value = "  -0.62%  "

$ws.Range($helper).Formula = "=""8.65"""
$ws.Range($helper).Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range($helper).Clear()
$ws.Range("E46").Value = "  +2.37%  "

$ws.Range($helper).Formula = "=""0.000281"""
$ws.Range($helper).Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range($helper).Clear()
$ws.Range("E48").Value = "  +23.07%  "

$ws.Range($helper).Formula = "=""0.0368"""
$ws.Range($helper).Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range($helper).Clear()
$ws.Range("E49").Value = "  +2.69%  "

$ws.Range("E50").Value = "  +0.14%  "

$ws.Range($helper).Formula = "=""40.18"""
$ws.Range($helper).Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range($helper).Clear()
$ws.Range("E51").Value = "  +1.21%  "

$excel.CutCopyMode = 0

